# Weekly update for the Coliflor (Mercado Mayorista Lo Valledor de Santiago)
# price sheet: a brand-new week (2021-09-09, serial 44448) is inserted as the
# earliest record, so the eight existing weekly rows (391-398, covering the
# Primera/Segunda pairs for 2020-12-02, 2021-02-11, 2021-07-22 and
# 2021-07-23) shift down two rows to 393-400, and the new week's data lands
# on rows 391-392.
#
# Columns: A Mercado ID, B Mercado, C Region, D Fecha (serial), E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 391; Fecha = 44448; Calidad = "Primera"; Volumen = 13300; Min = 400; Max = 500; Prom = 451 },
    @{ Row = 392; Fecha = 44448; Calidad = "Segunda"; Volumen = 3400;  Min = 300; Max = 300; Prom = 300 },
    @{ Row = 393; Fecha = 44167; Calidad = "Primera"; Volumen = 19200; Min = 600; Max = 700; Prom = 653 },
    @{ Row = 394; Fecha = 44167; Calidad = "Segunda"; Volumen = 7000;  Min = 500; Max = 500; Prom = 500 },
    @{ Row = 395; Fecha = 44238; Calidad = "Primera"; Volumen = 14000; Min = 800; Max = 900; Prom = 854 },
    @{ Row = 396; Fecha = 44238; Calidad = "Segunda"; Volumen = 5500;  Min = 700; Max = 700; Prom = 700 },
    @{ Row = 397; Fecha = 44399; Calidad = "Primera"; Volumen = 20100; Min = 500; Max = 600; Prom = 541 },
    @{ Row = 398; Fecha = 44399; Calidad = "Segunda"; Volumen = 6500;  Min = 400; Max = 500; Prom = 443 },
    @{ Row = 399; Fecha = 44400; Calidad = "Primera"; Volumen = 20200; Min = 500; Max = 600; Prom = 536 },
    @{ Row = 400; Fecha = 44400; Calidad = "Segunda"; Volumen = 6200;  Min = 400; Max = 500; Prom = 445 }
)

# The existing date column uses a custom date-time number format (applied
# via style index on rows that already existed); grab it once so the two
# brand-new rows (399, 400) get the same display format.
$dateFormat = $ws.Cells.Item(390, 4).NumberFormat

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"

    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 4).Value = $r.Fecha

    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112008
    $ws.Cells.Item($row, 7).Value = "Coliflor"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
